$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 19 Nov 2024 - combat implementation
$ws.Range("A3").Value = 45615
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B3").Value = "combat implementation - basic parts and placeholder sprites/abiklities"
$ws.Range("C3").Value = 6

# Row 4: 20 Nov 2024 - combat bugfixes
$ws.Range("A4").Value = 45616
$ws.Range("A4").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B4").Value = "combat bugfixes and extra features"
$ws.Range("C4").Value = 5

$ws.Range("D12").Select()
